# Updated symbol list on Fri Feb  3 09:52:27 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for the
# crypto rows on Sheet1. Values are written with a leading apostrophe so
# Excel stores them as literal text (matching the sheet's existing
# inline-string formatting) instead of auto-converting them to numbers
# or percentages; the style is then reset to "Normal" so no visible
# formatting/style change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'324.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.60%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-1.67%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.691"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'7.67%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.21%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.62%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.494"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.39%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.16%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.16%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9254"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1241"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-7.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'0.39%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'8.721"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09267"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.33%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.03629"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.65%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'9.49%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.001291"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.34%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006085"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-3.01%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.75%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.35%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'3.68%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2411"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-5.84%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04407"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.53%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001260"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.13%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004673"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'8.92%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001151"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-3.29%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02496"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.10%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05326"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.73%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007442"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.89%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009589"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.49%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1406"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.44%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002117"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.45%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01143"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'13.37%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006716"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.09%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.002972"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-11.13%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002290"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.60%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
